$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78.. down to 79..
$ws.Rows(78).Insert()

# Populate the newly inserted row 78 with the new weekly record
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 44880
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = 100112032
$ws.Cells.Item(78, 7).Value = "Zapallo italiano"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 100
$ws.Cells.Item(78, 11).Value = 10000
$ws.Cells.Item(78, 12).Value = 11000
$ws.Cells.Item(78, 13).Value = 10500
$ws.Cells.Item(78, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(78, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(78, 16).Value = 210
$ws.Cells.Item(78, 17).Value = 50
$ws.Cells.Item(78, 18).Value = "Hortaliza"
